$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the daily stats for rows 52-60 (columns B,C,E,F,G,H,I,J).
# Column D already holds the shared formula "=Bn-Cn" and recalculates
# automatically once B/C are populated.

$rowsData = @{
    52 = @{ B=3431; C=3085; E=87;  F=4;  G=5;  H=204; I=14; J=1796 }
    53 = @{ B=673;  C=108;  E=12;  F=1;  G=1;  H=32;  I=0;  J=0    }
    54 = @{ B=4932; C=4555; E=125; F=8;  G=9;  H=266; I=16; J=2624 }
    55 = @{ B=4661; C=4189; E=96;  F=3;  G=4;  H=234; I=3;  J=121  }
    56 = @{ B=3294; C=2785; E=84;  F=10; G=11; H=196; I=4;  J=705  }
    57 = @{ B=3623; C=3287; E=70;  F=5;  G=8;  H=241; I=9;  J=93   }
    58 = @{ B=3772; C=3278; E=76;  F=7;  G=7;  H=220; I=5;  J=408  }
    59 = @{ B=2879; C=2593; E=70;  F=2;  G=2;  H=109; I=2;  J=160  }
    60 = @{ B=424;  C=17;   E=19;              G=0;   H=17.9; I=0; J=0 }
}

foreach ($r in 52..60) {
    $vals = $rowsData[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}

# Update the frozen pane top-left cell and the active selection to match
# the scrolled-down view after the new rows were entered.
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H61").Select()
